# Daily attendance update - 2025-08-22
# Fills in the attendance grid (WFO/WFH/SL) for the week of 2025-08-18..2025-08-22
# (columns W..AA) on the "WCS_Team_August_2025" sheet, for every tracked employee
# (rows 3..18). The WFO/WFH summary counts in rows 21-36 are driven by COUNTIF
# formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WCS_Team_August_2025")
$ws.Activate()

# xlPasteFormats constant, used to copy the already-existing "Good"/"Neutral"
# cell styles (with their borders) from earlier same-row cells onto the new
# week's cells, instead of re-creating/duplicating style records.
$xlPasteFormats = -4122

function Set-AttendanceCell($TargetCell, $StyleSourceCell, $Text) {
    $ws.Range($StyleSourceCell).Copy() | Out-Null
    $ws.Range($TargetCell).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($TargetCell).Value = $Text
}

# row -> values for columns W, X, Y, Z, AA (2025-08-18 .. 2025-08-22)
$weekPlan = @{
    3  = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    4  = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    5  = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFO' }
    6  = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    7  = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    8  = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    9  = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    10 = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    11 = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    12 = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    13 = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    14 = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    15 = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    16 = @{ W='SL';  X='SL';  Y='SL';  Z='WFH'; AA='WFH' }
    17 = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
    18 = @{ W='WFO'; X='WFO'; Y='WFO'; Z='WFH'; AA='WFH' }
}

foreach ($row in 3..18) {
    $plan = $weekPlan[$row]
    if ($plan -eq $null) { continue }

    # donor cell for the WFH style is always column S of the same row.
    $wfhDonor = "S$row"
    # donor cell for the WFO/SL style: column I of the same row already
    # carries the right style for this employee that month (row 16 /
    # Swarnagowri is SL-styled there because she is on leave that month).
    $wfoDonor = "I$row"

    foreach ($col in @('W','X','Y','Z','AA')) {
        $text = $plan[$col]
        if ($text -eq 'WFH') {
            $donor = $wfhDonor
        } else {
            $donor = $wfoDonor
        }
        Set-AttendanceCell "$col$row" $donor $text
    }
}

# Recalculate so the WFO/WFH COUNTIF summary rows (21-36) pick up the new values.
$excel.Calculate()

# Restore the selection captured at save time.
$ws.Range("AA5").Select()
